$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.354.09"
$ws.Range("E2").Value = "  +3.78%  "

$ws.Range("D3").Value = "1.717.69"
$ws.Range("E3").Value = "  +3.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9985"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4710"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2640"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").Value = "1.710.81"
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07078"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.39%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5905"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.422"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "26.344.97"
$ws.Range("E18").Value = "  +3.80%  "

$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").Value = "1.930.70"
$ws.Range("E21").Value = "  +3.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.552"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.824"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.16%  "

$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.404"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.765"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.049"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.691"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07711"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04424"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.611"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6232"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9736"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9275"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "113.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.409"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.13%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.909"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.54%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.0000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.282"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.26%  "

$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1149"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.245"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05289"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.689"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.222"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("E51").Value = "  +1.18%  "
